# PTH 5mm axial capacitors.
# Fill in vendor (Digi-Key), vendor PN, manufacturer PN and unit price
# for the AXIAL-5MM capacitor BOM rows that were previously missing
# that sourcing data (rows 8, 9, 17, 20 and 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $VendorPN, $ManufacturerPN, $UnitPrice, $Note)

    $ws.Cells.Item($Row, 5).Value = "DK"
    $ws.Cells.Item($Row, 6).Value = $VendorPN
    $ws.Cells.Item($Row, 7).Value = $ManufacturerPN
    $ws.Cells.Item($Row, 9).Value = $UnitPrice

    if ($Note) {
        $ws.Cells.Item($Row, 12).Value = $Note
    }
}

Set-Row 8  "445-2904-ND" "FK22Y5V1E226Z" 1.25
Set-Row 9  "445-2880-ND" "FK26Y5V1E475Z" 0.43
Set-Row 17 "445-5257-ND" "FK28X7R1H683K" 0.35
Set-Row 20 "445-5253-ND" "FK28X7R1H153K" 0.29 "Real part value is 0.015uF"
Set-Row 21 "1292PH-ND"   "S680K29SL0N63J5R" 0.24

[void]$ws.Range("D22").Select()
